$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AA: the old AA (running total "STATUS"/label column)
# shifts to AB, and the new AA inherits column Z formatting (date style).
$ws.Columns.Item(27).Insert()

# Header row: the new AA1 gets the next day; AB1 keeps the shifted label cell as-is.
$ws.Range("AA1").Value = 44880

# Data rows 2-29: new AA is the updated running total for the new day,
# AB is the refreshed count that used to live in AA.
# Rows 30-34 previously held SUM formulas in AA (now shifted to AB); both
# columns are overwritten with literal numbers, clearing any formula.
$ws.Range("AA2").Value = 17035
$ws.Range("AB2").Value = 23
$ws.Range("AA3").Value = 1527
$ws.Range("AB3").Value = 6
$ws.Range("AA4").Value = 16211
$ws.Range("AB4").Value = 30
$ws.Range("AA5").Value = 16182
$ws.Range("AB5").Value = 20
$ws.Range("AA6").Value = 16128
$ws.Range("AB6").Value = 24
$ws.Range("AA7").Value = 16094
$ws.Range("AB7").Value = 17
$ws.Range("AA8").Value = 15888
$ws.Range("AB8").Value = 8
$ws.Range("AA9").Value = 15843
$ws.Range("AB9").Value = 13
$ws.Range("AA10").Value = 15762
$ws.Range("AB10").Value = 27
$ws.Range("AA11").Value = 15686
$ws.Range("AB11").Value = 17
$ws.Range("AA12").Value = 1457
$ws.Range("AB12").Value = 4
$ws.Range("AA13").Value = 15129
$ws.Range("AB13").Value = 8
$ws.Range("AA14").Value = 14003
$ws.Range("AB14").Value = 11
$ws.Range("AA15").Value = 13930
$ws.Range("AB15").Value = 6
$ws.Range("AA16").Value = 13871
$ws.Range("AB16").Value = 7
$ws.Range("AA17").Value = 12534
$ws.Range("AB17").Value = 5
$ws.Range("AA18").Value = 13414
$ws.Range("AB18").Value = 9
$ws.Range("AA19").Value = 1307
$ws.Range("AB19").Value = 3
$ws.Range("AA20").Value = 13303
$ws.Range("AB20").Value = 6
$ws.Range("AA21").Value = 13211
$ws.Range("AB21").Value = 1
$ws.Range("AA22").Value = 13165
$ws.Range("AB22").Value = 38
$ws.Range("AA23").Value = 1298
$ws.Range("AB23").Value = 1
$ws.Range("AA24").Value = 16570
$ws.Range("AB24").Value = 11
$ws.Range("AA25").Value = 16411
$ws.Range("AB25").Value = 5
$ws.Range("AA26").Value = 1511
$ws.Range("AB26").Value = 2
$ws.Range("AA27").Value = 16310
$ws.Range("AB27").Value = 19
$ws.Range("AA28").Value = 16124
$ws.Range("AB28").Value = 5
$ws.Range("AA29").Value = 16200
$ws.Range("AB29").Value = 4
$ws.Range("AA30").Value = 398
$ws.Range("AB30").Value = 0
$ws.Range("AA31").Value = 3534
$ws.Range("AB31").Value = 0
$ws.Range("AA32").Value = 3514
$ws.Range("AB32").Value = 0
$ws.Range("AA33").Value = 3504
$ws.Range("AB33").Value = 0
$ws.Range("AA34").Value = 3492
$ws.Range("AB34").Value = 0

# Refresh the view: clear the old scrolled-down position and select AB3
$ws.Range("AB3").Select()
